# NextTransaction.xaml deleted and replaced with actions in transitions.
# The "Workblocks" sheet had a dedicated wbNextTransaction_Type /
# wbNextTransaction_SuppressSuccessful pair of config rows (rows 11-12);
# since that workflow no longer exists as a separate state, remove those
# two rows and shift everything below them up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workblocks")

# Remove the wbNextTransaction_Type / wbNextTransaction_SuppressSuccessful
# rows; Excel shifts rows 13:18 up to 11:16.
$ws.Rows("11:12").Delete()

# Leave the selection where the author left it (end of the now-shorter
# list) and make Workblocks the active/visible tab.
$ws.Range("B18").Select() | Out-Null
$ws.Activate()
